$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2033.625
$ws.Range("J17").Value = 2033.625
$ws.Range("L17").Value = 6100.875
$ws.Range("N17").Value = -6436.875
$ws.Range("H28").Value = 2095.476
$ws.Range("I28").Value = 725.5714
$ws.Range("K28").Value = 725.5714
$ws.Range("M28").Value = -240.5714
$ws.Range("H64").Value = 8018.92
$ws.Range("I64").Value = 7998.909
$ws.Range("J64").Value = 8165.6665
$ws.Range("K64").Value = 7998.909
$ws.Range("L64").Value = 8165.6665
$ws.Range("M64").Value = -7750.909
$ws.Range("N64").Value = -8661.666499999999
$ws.Range("H67").Value = 8018.92
$ws.Range("I67").Value = 7998.909
$ws.Range("J67").Value = 8165.6665
$ws.Range("K67").Value = 7998.909
$ws.Range("L67").Value = 8165.6665
$ws.Range("M67").Value = -7140.909
$ws.Range("N67").Value = -9881.666499999999
$ws.Range("H86").Value = 2684
$ws.Range("I86").Value = 2227.25
$ws.Range("K86").Value = 2227.25
$ws.Range("M86").Value = -1104.25
$ws.Range("H89").Value = 2684
$ws.Range("I89").Value = 2227.25
$ws.Range("K89").Value = 11136.25
$ws.Range("M89").Value = -5520.25
$ws.Range("H105").Value = 33328
$ws.Range("J105").Value = 33328
$ws.Range("L105").Value = 33328
$ws.Range("N105").Value = -40316
$ws.Range("H132").Value = 1282.3182
$ws.Range("I132").Value = 1328.475
$ws.Range("J132").Value = 820.75
$ws.Range("K132").Value = 3985.425
$ws.Range("L132").Value = 2462.25
$ws.Range("M132").Value = -1455.425
$ws.Range("N132").Value = -7522.25
$ws.Range("H137").Value = 355723.06
$ws.Range("I137").Value = 1595.4849
$ws.Range("J137").Value = 1816499.2
$ws.Range("K137").Value = 4786.4547
$ws.Range("L137").Value = 5449497.6
$ws.Range("M137").Value = -2236.4547
$ws.Range("N137").Value = -5454597.6
$ws.Range("H141").Value = 2965.8235
$ws.Range("I141").Value = 2494.6
$ws.Range("K141").Value = 7483.799999999999
$ws.Range("M141").Value = -2303.799999999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4676.3076
$ws.Range("I32").Value = 1503.1228
$ws.Range("K32").Value = 1503.1228
$ws.Range("M32").Value = -1216.1228
$ws.Range("H61").Value = 50094.285
$ws.Range("I61").Value = 2424.125
$ws.Range("J61").Value = 202638.8
$ws.Range("K61").Value = 2424.125
$ws.Range("L61").Value = 202638.8
$ws.Range("M61").Value = -2212.125
$ws.Range("N61").Value = -203062.8
$ws.Range("H74").Value = 70150.13
$ws.Range("I74").Value = 93226.17999999999
$ws.Range("K74").Value = 93226.17999999999
$ws.Range("M74").Value = -92352.17999999999
$ws.Range("H77").Value = 70150.13
$ws.Range("I77").Value = 93226.17999999999
$ws.Range("K77").Value = 466130.9
$ws.Range("M77").Value = -461762.9
$ws.Range("H136").Value = 50094.285
$ws.Range("I136").Value = 2424.125
$ws.Range("J136").Value = 202638.8
$ws.Range("K136").Value = 7272.375
$ws.Range("L136").Value = 607916.3999999999
$ws.Range("M136").Value = -4722.375
$ws.Range("N136").Value = -613016.3999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 5088.769
$ws.Range("I86").Value = 5474
$ws.Range("J86").Value = 4222
$ws.Range("K86").Value = 5474
$ws.Range("L86").Value = 4222
$ws.Range("M86").Value = -4351
$ws.Range("N86").Value = -6468
$ws.Range("H89").Value = 5088.769
$ws.Range("I89").Value = 5474
$ws.Range("J89").Value = 4222
$ws.Range("K89").Value = 27370
$ws.Range("L89").Value = 21110
$ws.Range("M89").Value = -21754
$ws.Range("N89").Value = -32342
$ws.Range("H105").Value = 33165.03
$ws.Range("I105").Value = 36199.38
$ws.Range("J105").Value = 3833
$ws.Range("K105").Value = 36199.38
$ws.Range("L105").Value = 3833
$ws.Range("M105").Value = -34452.38
$ws.Range("N105").Value = -7327
$ws.Range("H107").Value = 9095496
$ws.Range("I107").Value = 12504294
$ws.Range("J107").Value = 5369.6665
$ws.Range("K107").Value = 12504294
$ws.Range("L107").Value = 5369.6665
$ws.Range("M107").Value = -12502374
$ws.Range("N107").Value = -9209.666499999999
$ws.Range("H134").Value = 5100.303
$ws.Range("I134").Value = 2309.7727
$ws.Range("J134").Value = 10681.363
$ws.Range("K134").Value = 6929.3181
$ws.Range("L134").Value = 32044.089
$ws.Range("M134").Value = -4394.3181
$ws.Range("N134").Value = -37114.089

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3594.5789
$ws.Range("I31").Value = 2209.625
$ws.Range("J31").Value = 4601.8184
$ws.Range("K31").Value = 2209.625
$ws.Range("L31").Value = 4601.8184
$ws.Range("M31").Value = -1914.625
$ws.Range("N31").Value = -5191.8184
$ws.Range("H34").Value = 3594.5789
$ws.Range("I34").Value = 2209.625
$ws.Range("J34").Value = 4601.8184
$ws.Range("K34").Value = 2209.625
$ws.Range("L34").Value = 4601.8184
$ws.Range("M34").Value = -2007.625
$ws.Range("N34").Value = -5005.8184
$ws.Range("H99").Value = 11865227
$ws.Range("I99").Value = 13890506
$ws.Range("J99").Value = 7814669
$ws.Range("K99").Value = 13890506
$ws.Range("L99").Value = 7814669
$ws.Range("M99").Value = -13889008
$ws.Range("N99").Value = -7817665
$ws.Range("H126").Value = 11865227
$ws.Range("I126").Value = 13890506
$ws.Range("J126").Value = 7814669
$ws.Range("K126").Value = 41671518
$ws.Range("L126").Value = 23444007
$ws.Range("M126").Value = -41669048
$ws.Range("N126").Value = -23448947

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 73204.8
$ws.Range("J37").Value = 73204.8
$ws.Range("L37").Value = 219614.4
$ws.Range("N37").Value = -219838.4
$ws.Range("H55").Value = 818.1818
$ws.Range("I55").Value = 818.1818
$ws.Range("K55").Value = 2454.5454
$ws.Range("M55").Value = -2277.5454
$ws.Range("H56").Value = 6594.143
$ws.Range("I56").Value = 6594.143
$ws.Range("K56").Value = 6594.143
$ws.Range("M56").Value = -6064.143
$ws.Range("H107").Value = 1963.6666
$ws.Range("J107").Value = 1995
$ws.Range("L107").Value = 5985
$ws.Range("N107").Value = -9825
$ws.Range("H119").Value = 745.5
$ws.Range("I119").Value = 316.66666
$ws.Range("J119").Value = 2032
$ws.Range("K119").Value = 949.9999799999999
$ws.Range("L119").Value = 6096
$ws.Range("M119").Value = 3888.00002
$ws.Range("N119").Value = -15772
$ws.Range("H131").Value = 1181.1111
$ws.Range("I131").Value = 906.6667
$ws.Range("J131").Value = 1730
$ws.Range("K131").Value = 2720.0001
$ws.Range("L131").Value = 5190
$ws.Range("M131").Value = 2319.9999
$ws.Range("N131").Value = -15270

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 9161.5
$ws.Range("I70").Value = 9593.799999999999
$ws.Range("K70").Value = 9593.799999999999
$ws.Range("M70").Value = -9323.799999999999
$ws.Range("H73").Value = 9161.5
$ws.Range("I73").Value = 9593.799999999999
$ws.Range("K73").Value = 9593.799999999999
$ws.Range("M73").Value = -8657.799999999999
$ws.Range("H97").Value = 477.58334
$ws.Range("I97").Value = 324.33334
$ws.Range("J97").Value = 937.3333
$ws.Range("K97").Value = 324.33334
$ws.Range("L97").Value = 937.3333
$ws.Range("M97").Value = 171.66666
$ws.Range("N97").Value = -1929.3333
$ws.Range("H126").Value = 3778.7693
$ws.Range("I126").Value = 2440
$ws.Range("J126").Value = 4615.5
$ws.Range("K126").Value = 7320
$ws.Range("L126").Value = 13846.5
$ws.Range("M126").Value = -4850
$ws.Range("N126").Value = -18786.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 2222920.2
$ws.Range("I55").Value = 622.63635
$ws.Range("K55").Value = 622.63635
$ws.Range("M55").Value = -449.63635
$ws.Range("H93").Value = 2070.3845
$ws.Range("I93").Value = 1658.6666
$ws.Range("K93").Value = 1658.6666
$ws.Range("M93").Value = -410.6666
$ws.Range("H132").Value = 2826.5186
$ws.Range("I132").Value = 2564.6316
$ws.Range("K132").Value = 7693.8948
$ws.Range("M132").Value = -5163.8948
$ws.Range("H136").Value = 1674.7273
$ws.Range("I136").Value = 1587.8096
$ws.Range("J136").Value = 3500
$ws.Range("K136").Value = 4763.4288
$ws.Range("L136").Value = 10500
$ws.Range("M136").Value = -2213.4288
$ws.Range("N136").Value = -15600

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1652.6111
$ws.Range("I122").Value = 1249.9
$ws.Range("J122").Value = 2156
$ws.Range("K122").Value = 3749.7
$ws.Range("L122").Value = 6468
$ws.Range("M122").Value = -1299.7
$ws.Range("N122").Value = -11368
